# Refresh the cryptos list (price + 1h volume change %) in Sheet1.
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# A handful of rows also swap two adjacent coins' rank order (name/link/price/volume).
# Cells whose new "Price" text parses as a plain number (e.g. "1.00", "0.610")
# are forced to Text format first so Excel doesn't collapse them into a numeric
# value (which would drop the trailing zeros / formatting shown in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.475.26"
$ws.Cells.Item(2, 5).Value = "  +3.30%  "
$ws.Cells.Item(3, 4).Value = "2.646.03"
$ws.Cells.Item(3, 5).Value = "  +1.20%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.27%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "567.05"
$ws.Cells.Item(5, 5).Value = "  +6.15%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "147.11"
$ws.Cells.Item(6, 5).Value = "  +3.03%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.997"
$ws.Cells.Item(7, 5).Value = "  -0.26%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.610"
$ws.Cells.Item(8, 5).Value = "  +6.52%  "
$ws.Cells.Item(9, 4).Value = "2.673.24"
$ws.Cells.Item(9, 5).Value = "  +2.03%  "
$ws.Cells.Item(10, 5).Value = "  -0.28%  "
$ws.Cells.Item(11, 5).Value = "  +5.21%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.145"
$ws.Cells.Item(12, 5).Value = "  +7.17%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.343"
$ws.Cells.Item(13, 5).Value = "  +3.11%  "
$ws.Cells.Item(14, 4).Value = "3.114.94"
$ws.Cells.Item(14, 5).Value = "  +0.96%  "
$ws.Cells.Item(15, 4).Value = "60.444.46"
$ws.Cells.Item(15, 5).Value = "  +3.38%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "21.91"
$ws.Cells.Item(16, 5).Value = "  +5.73%  "
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.0000137"
$ws.Cells.Item(17, 5).Value = "  +4.45%  "
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "2.653.01"
$ws.Cells.Item(18, 5).Value = "  +1.22%  "
$ws.Cells.Item(19, 5).Value = "  +3.55%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "344.27"
$ws.Cells.Item(20, 5).Value = "  +3.14%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.48"
$ws.Cells.Item(21, 5).Value = "  +3.48%  "
$ws.Cells.Item(22, 5).Value = "  +2.51%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.83"
$ws.Cells.Item(23, 5).Value = "  +1.38%  "
$ws.Cells.Item(24, 5).Value = "  -0.11%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "66.70"
$ws.Cells.Item(25, 5).Value = "  +0.51%  "
$ws.Cells.Item(26, 5).Value = "  +5.37%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.164"
$ws.Cells.Item(27, 5).Value = "  +1.48%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.994"
$ws.Cells.Item(28, 5).Value = "  -0.51%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.40"
$ws.Cells.Item(29, 5).Value = "  +4.40%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0782"
$ws.Cells.Item(30, 5).Value = "  +6.71%  "
$ws.Cells.Item(31, 2).Value = "USDe"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.998"
$ws.Cells.Item(31, 5).Value = "  -0.08%  "
$ws.Cells.Item(32, 2).Value = "Aptos"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.32"
$ws.Cells.Item(32, 5).Value = "  +7.65%  "
$ws.Cells.Item(33, 5).Value = "  +4.80%  "
$ws.Cells.Item(34, 2).Value = "Monero"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "155.11"
$ws.Cells.Item(34, 5).Value = "  +1.83%  "
$ws.Cells.Item(35, 2).Value = "EthereumClassic"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "19.25"
$ws.Cells.Item(35, 5).Value = "  +2.09%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.10"
$ws.Cells.Item(36, 5).Value = "  +5.34%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.906"
$ws.Cells.Item(37, 5).Value = "  +7.35%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.909"
$ws.Cells.Item(38, 5).Value = "  +11.69%  "
$ws.Cells.Item(39, 5).Value = "  +5.92%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "37.58"
$ws.Cells.Item(40, 5).Value = "  +1.18%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.52"
$ws.Cells.Item(41, 5).Value = "  +7.14%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "302.66"
$ws.Cells.Item(42, 5).Value = "  +7.44%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.67"
$ws.Cells.Item(43, 5).Value = "  +2.77%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.994"
$ws.Cells.Item(44, 5).Value = "  -0.61%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0985"
$ws.Cells.Item(45, 5).Value = "  +4.72%  "
$ws.Cells.Item(46, 5).Value = "  +2.23%  "
$ws.Cells.Item(47, 5).Value = "  +4.09%  "
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "128.53"
$ws.Cells.Item(48, 5).Value = "  +12.81%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "19.46"
$ws.Cells.Item(49, 5).Value = "  +2.51%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "10.69"
$ws.Cells.Item(50, 5).Value = "  -0.03%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0236"
$ws.Cells.Item(51, 5).Value = "  +5.18%  "
